$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.690.97'
$ws.Range("E2").Value = '  -1.30%  '

$ws.Range("D3").Value = '3.074.57'
$ws.Range("E3").Value = '  -2.64%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.08'
$ws.Range("E5").Value = '  +9.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.05'
$ws.Range("E6").Value = '  -1.44%  '

$ws.Range("E7").Value = '  -10.45%  '

$ws.Range("E8").Value = '  -1.85%  '

$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").Value = '3.075.27'
$ws.Range("E10").Value = '  -2.56%  '

$ws.Range("E11").Value = '  -6.21%  '

$ws.Range("E12").Value = '  -1.70%  '

$ws.Range("E13").Value = '  +1.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.31'
$ws.Range("E14").Value = '  +0.52%  '

$ws.Range("D15").Value = '89.703.46'
$ws.Range("E15").Value = '  -1.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.37'
$ws.Range("E16").Value = '  -6.77%  '

$ws.Range("D17").Value = '3.669.81'
$ws.Range("E17").Value = '  -1.96%  '

$ws.Range("D18").Value = '3.091.26'
$ws.Range("E18").Value = '  -1.13%  '

$ws.Range("E19").Value = '  +2.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000212'
$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.79'
$ws.Range("E21").Value = '  -5.84%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '433.47'
$ws.Range("E22").Value = '  -8.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.41'
$ws.Range("E23").Value = '  +3.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.78'
$ws.Range("E24").Value = '  -4.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.59'
$ws.Range("E25").Value = '  -5.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '86.64'
$ws.Range("E26").Value = '  -9.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.76'
$ws.Range("E27").Value = '  -4.82%  '

$ws.Range("E28").Value = '  -1.65%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.06'
$ws.Range("E30").Value = '  -3.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("E32").Value = '  -3.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.195'
$ws.Range("E33").Value = '  -8.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.152'
$ws.Range("E34").Value = '  +4.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.57'
$ws.Range("E35").Value = '  -7.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.72'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.12'
$ws.Range("E37").Value = '  +2.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '495.69'
$ws.Range("E38").Value = '  -4.64%  '

$ws.Range("E39").Value = '  -3.14%  '

$ws.Range("E40").Value = '  -3.78%  '

$ws.Range("E41").Value = '  -3.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.59'
$ws.Range("E42").Value = '  +53.86%  '

$ws.Range("E43").Value = '  -0.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.396'
$ws.Range("E45").Value = '  -6.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.58'
$ws.Range("E46").Value = '  +0.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.85'
$ws.Range("E47").Value = '  -6.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.675'
$ws.Range("E48").Value = '  -9.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.34'
$ws.Range("E49").Value = '  -2.73%  '

$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.30'
$ws.Range("E51").Value = '  -4.25%  '
